$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from existing template rows onto the new rows ---
# Rows 8,9,13,14 follow the same look as row 2 (full A:E, style 4/5/5/5/5)
$ws.Range("A2:E2").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)

# Rows 10,11 follow row 5's look (no A cell, style 4/5/5/5 on B:E)
$ws.Range("B5:E5").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)
$ws.Range("B5:E5").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)

# Row 12 follows row 6's look (style 6/6/7/7/7)
$ws.Range("A6:E6").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 8 ---
$ws.Cells.Item(8, 1).Value = "SCRIPT/G01P03A/um1310.ssb"
$ws.Cells.Item(8, 2).Value = 238
$ws.Cells.Item(8, 3).Value = " I\'m interested in a treasure\ncalled [CS:I]Mystery Part[CR]."
$ws.Cells.Item(8, 4).Value = " Я интересуюсь сокровищем,\nназываемом [CS:I]Загадочной Частью[CR]."
$ws.Cells.Item(8, 5).Value = " Ÿ éîóåñåòôýòû òïëñïâéþåí,\nîàèúâàåíïí [CS:I]Èàãàäïœîïê Œàòóûý[CR]."

# --- Row 9 ---
$ws.Cells.Item(9, 1).Value = "SCRIPT/G01P03A/um1402.ssb"
$ws.Cells.Item(9, 2).Value = 241
$ws.Cells.Item(9, 3).Value = " It\'s a legendary item that\'s\nveiled in mystery.[K] Well? Don\'t you find it\nintriguing?"
$ws.Cells.Item(9, 4).Value = " Это окутанный тайной легендарный\nпредмет.[K] Ну как? Разве тебе не интересно?"
$ws.Cells.Item(9, 5).Value = " Üóï ïëôóàîîúê óàêîïê ìåãåîäàñîúê\nðñåäíåó.[K] Îô ëàë? Ñàèâå óåáå îå éîóåñîï?"

# --- Row 10 ---
$ws.Cells.Item(10, 2).Value = 244
$ws.Cells.Item(10, 3).Value = " But my partner [CS:N]Murkrow[CR]…"
$ws.Cells.Item(10, 4).Value = " Но вот моя спутница [CS:N]Маркроу[CR]..."
$ws.Cells.Item(10, 5).Value = " Îï âïó íïÿ òðôóîéøà [CS:N]Íàñëñïô[CR]..."

# --- Row 11 ---
$ws.Cells.Item(11, 2).Value = 247
$ws.Cells.Item(11, 3).Value = " She\'s like, \`"I have no interest in\na drab and dingy thing like that!\`""
$ws.Cells.Item(11, 4).Value = " Она такая: \`"Меня не интересует\nтакая серая и тусклая вещь!\`""
$ws.Cells.Item(11, 5).Value = " Ïîà óàëàÿ: \`"Íåîÿ îå éîóåñåòôåó\nóàëàÿ òåñàÿ é óôòëìàÿ âåþû!\`""

# --- Row 12 ---
$ws.Cells.Item(12, 2).Value = 250
$ws.Cells.Item(12, 3).Value = " Oh, what am I to do with her...?"
$ws.Cells.Item(12, 4).Value = " О, ну что мне с ней делать?.."
$ws.Cells.Item(12, 5).Value = " Ï, îô œóï íîå ò îåê äåìàóû?.."

# --- Row 13 ---
$ws.Cells.Item(13, 1).Value = "SCRIPT/T01P02A/um1410.ssb"
$ws.Cells.Item(13, 2).Value = 216
$ws.Cells.Item(13, 3).Value = " We found a mysterious\ntreasure box in a dungeon."
$ws.Cells.Item(13, 4).Value = " В подземелье нам попалась\nзагадочная шкатулка."
$ws.Cells.Item(13, 5).Value = " Â ðïäèåíåìûå îàí ðïðàìàòû\nèàãàäïœîàÿ šëàóôìëà."

# --- Row 14 ---
$ws.Cells.Item(14, 1).Value = "SCRIPT/T01P02A/um1505.ssb"
$ws.Cells.Item(14, 2).Value = 219
$ws.Cells.Item(14, 3).Value = " We visited Treasure Town so\nwe could have the box appraised at\n[CS:K]Xatu[CR]\'s shop."
$ws.Cells.Item(14, 4).Value = " Мы посетили Город Сокровищ,\nчтобы Оценщик [CS:K]Ксату[CR] открыл её нам."
$ws.Cells.Item(14, 5).Value = " Íú ðïòåóéìé Ãïñïä Òïëñïâéþ,\nœóïáú Ïøåîþéë [CS:K]Ëòàóô[CR] ïóëñúì åæ îàí."

# --- Row heights (10 and 12 stay at the sheet default, so no explicit set) ---
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 31.8
$ws.Rows.Item(13).RowHeight = 43.2
$ws.Rows.Item(14).RowHeight = 43.2

# --- View state ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D16").Select()
